$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-11 (years 2000-2009), shifting remaining rows (2010,2011,2012) up.
$ws.Range("A2:B11").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
